$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Item_Name: Palak -> Mint leaves
$ws.Range("A2").Value = "Mint leaves"

# Item_Code: PL_0001 -> Mint_0010
$ws.Range("D2").Value = "Mint_0010"

# UnitPrice: 20.56 -> 10.2
$ws.Range("E2").Value = 10.2

# MinStock: 1000 -> 15
$ws.Range("G2").Value = 15

# Description: (empty) -> Testing
$ws.Range("H2").Value = "Testing"
